# Auto-generated edit script applying the diff to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.560.24'
$ws.Range("E2").Value = '''  +3.66%  '
$ws.Range("D3").Value = '''1.694.01'
$ws.Range("E3").Value = '''  +2.22%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '''  +0.13%  '
$ws.Range("D5").Value = '''317.27'
$ws.Range("E5").Value = '''  +2.42%  '
$ws.Range("E6").Value = '''  +0.07%  '
$ws.Range("D7").Value = '''0.3951'
$ws.Range("E7").Value = '''  +2.04%  '
$ws.Range("D8").Value = '''0.4012'
$ws.Range("E8").Value = '''  +1.84%  '
$ws.Range("D9").Value = '''1.535'
$ws.Range("E9").Value = '''  +7.13%  '
$ws.Range("D10").Value = '''54.30'
$ws.Range("E10").Value = '''  +10.70%  '
$ws.Range("D11").Value = '''1.002'
$ws.Range("E11").Value = '''  +0.12%  '
$ws.Range("D12").Value = '''0.08762'
$ws.Range("E12").Value = '''  +1.50%  '
$ws.Range("E13").Value = '''  +8.57%  '
$ws.Range("D14").Value = '''23.24'
$ws.Range("E14").Value = '''  +2.96%  '
$ws.Range("D15").Value = '''0.00001322'
$ws.Range("E15").Value = '''  +0.69%  '
$ws.Range("D16").Value = '''7.582'
$ws.Range("E16").Value = '''  +4.88%  '
$ws.Range("D17").Value = '''1.697.04'
$ws.Range("D18").Value = '''100.86'
$ws.Range("E18").Value = '''  +1.31%  '
$ws.Range("D19").Value = '''0.07028'
$ws.Range("E19").Value = '''  +3.63%  '
$ws.Range("D20").Value = '''19.61'
$ws.Range("E20").Value = '''  +3.13%  '
$ws.Range("D21").Value = '''6.852'
$ws.Range("E21").Value = '''  +2.99%  '
$ws.Range("E22").Value = '''  +0.10%  '
$ws.Range("D23").Value = '''14.03'
$ws.Range("E23").Value = '''  +1.10%  '
$ws.Range("D24").Value = '''24.549.82'
$ws.Range("E24").Value = '''  +3.66%  '
$ws.Range("D25").Value = '''3.019'
$ws.Range("E25").Value = '''  +7.77%  '
$ws.Range("D26").Value = '''2.317'
$ws.Range("E26").Value = '''  -0.13%  '
$ws.Range("D27").Value = '''22.30'
$ws.Range("E27").Value = '''  +2.75%  '
$ws.Range("D28").Value = '''159.24'
$ws.Range("E28").Value = '''  +0.46%  '
$ws.Range("D29").Value = '''5.198'
$ws.Range("E29").Value = '''  +1.32%  '
$ws.Range("D30").Value = '''134.07'
$ws.Range("E30").Value = '''  +3.61%  '
$ws.Range("D31").Value = '''7.548'
$ws.Range("E31").Value = '''  +17.61%  '
$ws.Range("D32").Value = '''1.885.21'
$ws.Range("E32").Value = '''  +2.41%  '
$ws.Range("D33").Value = '''1.095'
$ws.Range("E33").Value = '''  -3.16%  '
$ws.Range("D34").Value = '''7.310'
$ws.Range("E34").Value = '''  +12.72%  '
$ws.Range("D35").Value = '''0.08539'
$ws.Range("E35").Value = '''  -0.01%  '
$ws.Range("D36").Value = '''11.36'
$ws.Range("E36").Value = '''  +9.67%  '
$ws.Range("D37").Value = '''1.974'
$ws.Range("E37").Value = '''  +0.57%  '
$ws.Range("D38").Value = '''0.2724'
$ws.Range("E38").Value = '''  +3.36%  '
$ws.Range("D39").Value = '''14.57'
$ws.Range("E39").Value = '''  +1.21%  '
$ws.Range("D40").Value = '''0.02752'
$ws.Range("E40").Value = '''  +9.19%  '
$ws.Range("D41").Value = '''0.09044'
$ws.Range("E41").Value = '''  +3.17%  '
$ws.Range("D42").Value = '''1.466'
$ws.Range("D43").Value = '''0.7675'
$ws.Range("E43").Value = '''  +1.95%  '
$ws.Range("D44").Value = '''0.7193'
$ws.Range("E44").Value = '''  +2.82%  '
$ws.Range("D45").Value = '''15.37'
$ws.Range("E45").Value = '''  +3.13%  '
$ws.Range("D46").Value = '''2.508'
$ws.Range("E46").Value = '''  +4.05%  '
$ws.Range("D47").Value = '''4.220'
$ws.Range("E47").Value = '''  +3.11%  '
$ws.Range("B48").Value = 'Flow'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D48").Value = '''1.353'
$ws.Range("E48").Value = '''  +14.81%  '
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '''1.001'
$ws.Range("E49").Value = '''  +0.09%  '
$ws.Range("D50").Value = '''141.29'
$ws.Range("E50").Value = '''  +2.49%  '
$ws.Range("D51").Value = '''0.08029'
$ws.Range("E51").Value = '''  +3.09%  '
